# Applies the "xlsx header rename / temp column consolidation" edit:
#  - Row 1 header labels: "t" -> "temp_lab", "mg"/"cl"/"br"/"na" -> "Mg"/"Cl"/"Br"/"Na"
#  - Row 3 recording date corrected (31959 -> 33420), matching the rest of the column
#  - Selection changed to the full column A (sqref A1:A1048576)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (column headers) -------------------------------------------------
$ws.Range("F1").Value = "temp_lab"
$ws.Range("J1").Value = "Mg"
$ws.Range("K1").Value = "Cl"
$ws.Range("L1").Value = "Br"
$ws.Range("M1").Value = "Na"

# --- Row 3 data fix ----------------------------------------------------------
$ws.Range("E3").Value = 33420

# --- Selection ---------------------------------------------------------------
$ws.Columns("A:A").Select()
